$wb = $excel.ActiveWorkbook

# Remove the extra sheets (Sheet2, Sheet3), keeping only the first sheet
$excel.DisplayAlerts = $false
for ($i = $wb.Worksheets.Count; $i -ge 2; $i--) {
    $wb.Worksheets.Item($i).Delete()
}

# Rename the remaining sheet
$ws = $wb.Worksheets.Item(1)
$ws.Name = "ValidLogin"

# Clear old contents and write the login form layout
$ws.Cells.Clear()
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "ADMIN"
$ws.Range("B2").Value = "manager"

# Select A3 to match the saved selection state
$ws.Range("A3").Select()

$wb.Save()
